$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BH (60): "Agosto.2021" report ------------------------------

# Header cell: same bold/centered/bordered look as the rest of row 1 (copy
# format from the previous header cell, then set the text).
$ws.Cells.Item(1, 60).Value = "Agosto.2021"
$ws.Cells.Item(1, 59).Copy()
$ws.Cells.Item(1, 60).PasteSpecial(-4122)  # xlPasteFormats

# Rows 2..73: column BH repeats the value already published in column BG (59)
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 60).Value = $ws.Cells.Item($r, 59).Value2
}

# Row 74: the new report revises the last published figure
$ws.Cells.Item(74, 60).Value = 29624

# --- New row 75: quarter 01-04-2021, first published in this report --------

# Column A holds a plain text label (not a date!) just like every other cell
# in that column, so force text with a quote-prefix, then drop the
# quote-prefix formatting it implies (copy the plain format from A74 back
# on top) so the cell ends up with the same "no explicit style" look as its
# neighbours above.
$ws.Cells.Item(75, 1).Value = "'01-04-2021"
$ws.Cells.Item(74, 1).Copy()
$ws.Cells.Item(75, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(75, 60).Value = 32833
